$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells are stored as text in the source data; the values
# look numeric, so force a text format before/while writing and then restore
# the default ("Normal") style so no stray number-format style is left behind.
$priceUpdates = @{
    "D2" = "243.67"
    "D3" = "23.13"
    "D4" = "5.420"
    "D6" = "3.455"
    "D7" = "6.531"
    "D8" = "0.8131"
    "D9" = "0.9139"
    "D10" = "0.1408"
    "D11" = "0.07466"
    "D12" = "0.03292"
    "D13" = "0.03056"
    "D14" = "0.09356"
    "D15" = "3.848"
    "D16" = "0.001582"
    "D17" = "0.04672"
    "D18" = "0.01123"
    "D19" = "0.006116"
    "D20" = "0.004996"
    "D21" = "0.0009819"
    "D23" = "3.608"
    "D24" = "2.136"
    "D40" = "0.03939"
    "D41" = "0.006206"
    "D42" = "0.1075"
    "D43" = "0.003001"
    "D44" = "0.007822"
    "D45" = "0.00005239"
    "D48" = "0.8984"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Volume(1h) (column E) cells are plain text; no special formatting required.
$volumeUpdates = @{
    "E18" = "17OneONEBestin24h"
    "E22" = "21NitroExNTX"
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
